$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a handful of neighbourhood names in column A to match the
# updated dataset's spelling/punctuation.
$ws.Range("A41").Value = "St.Andrew-Windfields"
$ws.Range("A46").Value = "Parkwoods-Donalda"
$ws.Range("A57").Value = "Leaside-Bennington"
$ws.Range("A60").Value = "Danforth East York"
$ws.Range("A62").Value = "Taylor-Massey"
$ws.Range("A72").Value = "Cabbagetown-South St.James Town"
$ws.Range("A75").Value = "North St.James Town"
$ws.Range("A92").Value = "Weston-Pellam Park"
$ws.Range("A120").Value = "Wexford/Maryvale"

# Update the view scroll position / active selection left over from the
# editing session.
$ws.Range("C55").Select()
$excel.ActiveWindow.ScrollRow = 39
